# Auto-generated edit script applying scheduled market-data refresh to Midgardsormr_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4990.364
$ws.Range("J40").Value = 4700.25
$ws.Range("L40").Value = 4700.25
$ws.Range("N40").Value = -5050.25
$ws.Range("H42").Value = 312.36365
$ws.Range("I42").Value = 340.6
$ws.Range("J42").Value = 30
$ws.Range("K42").Value = 1021.8
$ws.Range("L42").Value = 90
$ws.Range("M42").Value = -791.8000000000001
$ws.Range("N42").Value = -550
$ws.Range("H64").Value = 34708.43
$ws.Range("J64").Value = 9500
$ws.Range("L64").Value = 9500
$ws.Range("N64").Value = -9996
$ws.Range("H67").Value = 34708.43
$ws.Range("J67").Value = 9500
$ws.Range("L67").Value = 9500
$ws.Range("N67").Value = -11216
$ws.Range("H70").Value = 54407
$ws.Range("I70").Value = 6939.4
$ws.Range("J70").Value = 97559.37
$ws.Range("K70").Value = 20818.2
$ws.Range("L70").Value = 292678.11
$ws.Range("M70").Value = -20548.2
$ws.Range("N70").Value = -293218.11
$ws.Range("H73").Value = 54407
$ws.Range("I73").Value = 6939.4
$ws.Range("J73").Value = 97559.37
$ws.Range("K73").Value = 20818.2
$ws.Range("L73").Value = 292678.11
$ws.Range("M73").Value = -19882.2
$ws.Range("N73").Value = -294550.11
$ws.Range("H86").Value = 1901.125
$ws.Range("I86").Value = 1818.6364
$ws.Range("J86").Value = 2082.6
$ws.Range("K86").Value = 1818.6364
$ws.Range("L86").Value = 2082.6
$ws.Range("M86").Value = -695.6364000000001
$ws.Range("N86").Value = -4328.6
$ws.Range("H89").Value = 1901.125
$ws.Range("I89").Value = 1818.6364
$ws.Range("J89").Value = 2082.6
$ws.Range("K89").Value = 9093.182000000001
$ws.Range("L89").Value = 10413
$ws.Range("M89").Value = -3477.182000000001
$ws.Range("N89").Value = -21645
$ws.Range("H116").Value = 16624.033
$ws.Range("I116").Value = 16329.5
$ws.Range("K116").Value = 16329.5
$ws.Range("M116").Value = -12887.5
$ws.Range("H141").Value = 1302
$ws.Range("I141").Value = 1341.7142
$ws.Range("K141").Value = 4025.1426
$ws.Range("M141").Value = 1154.8574

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14200.26
$ws.Range("I32").Value = 15037.319
$ws.Range("K32").Value = 15037.319
$ws.Range("M32").Value = -14750.319
$ws.Range("H45").Value = 5304.3335
$ws.Range("I45").Value = 3992.25
$ws.Range("K45").Value = 3992.25
$ws.Range("M45").Value = -3615.25
$ws.Range("H61").Value = 7809.95
$ws.Range("I61").Value = 1500
$ws.Range("J61").Value = 11207.615
$ws.Range("K61").Value = 1500
$ws.Range("L61").Value = 11207.615
$ws.Range("M61").Value = -1288
$ws.Range("N61").Value = -11631.615
$ws.Range("H88").Value = 4385.0625
$ws.Range("J88").Value = 5763.273
$ws.Range("L88").Value = 5763.273
$ws.Range("N88").Value = -6575.273
$ws.Range("H91").Value = 4385.0625
$ws.Range("J91").Value = 5763.273
$ws.Range("L91").Value = 5763.273
$ws.Range("N91").Value = -8571.273000000001
$ws.Range("H136").Value = 7809.95
$ws.Range("I136").Value = 1500
$ws.Range("J136").Value = 11207.615
$ws.Range("K136").Value = 4500
$ws.Range("L136").Value = 33622.845
$ws.Range("M136").Value = -1950
$ws.Range("N136").Value = -38722.845

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 14481.261
$ws.Range("I20").Value = 23199.785
$ws.Range("K20").Value = 23199.785
$ws.Range("M20").Value = -22952.785
$ws.Range("H99").Value = 3960.1428
$ws.Range("I99").Value = 3850
$ws.Range("K99").Value = 3850
$ws.Range("M99").Value = -2352
$ws.Range("H105").Value = 3926.818
$ws.Range("I105").Value = 3483.1667
$ws.Range("K105").Value = 3483.1667
$ws.Range("M105").Value = -1736.1667
$ws.Range("H107").Value = 6125.5835
$ws.Range("J107").Value = 6332
$ws.Range("L107").Value = 6332
$ws.Range("N107").Value = -10172
$ws.Range("H134").Value = 1760.4546
$ws.Range("I134").Value = 1073.8235
$ws.Range("K134").Value = 3221.4705
$ws.Range("M134").Value = -686.4704999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 6681.8
$ws.Range("J2").Value = 5
$ws.Range("L2").Value = 5
$ws.Range("N2").Value = -231
$ws.Range("H8").Value = 199
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("H15").Value = 513.3333
$ws.Range("J15").Value = 545.5
$ws.Range("L15").Value = 545.5
$ws.Range("N15").Value = -885.5
$ws.Range("H22").Value = 731.9091
$ws.Range("I22").Value = 540.2
$ws.Range("J22").Value = 891.6667
$ws.Range("K22").Value = 540.2
$ws.Range("L22").Value = 891.6667
$ws.Range("M22").Value = -190.2
$ws.Range("N22").Value = -1591.6667
$ws.Range("H62").Value = 4413.857
$ws.Range("I62").Value = 4494.909
$ws.Range("J62").Value = 4116.6665
$ws.Range("K62").Value = 4494.909
$ws.Range("L62").Value = 4116.6665
$ws.Range("M62").Value = -3870.909
$ws.Range("N62").Value = -5364.6665
$ws.Range("H65").Value = 4413.857
$ws.Range("I65").Value = 4494.909
$ws.Range("J65").Value = 4116.6665
$ws.Range("K65").Value = 22474.545
$ws.Range("L65").Value = 20583.3325
$ws.Range("M65").Value = -19354.545
$ws.Range("N65").Value = -26823.3325
$ws.Range("H105").Value = 1068.3889
$ws.Range("I105").Value = 795.6875
$ws.Range("K105").Value = 795.6875
$ws.Range("M105").Value = 951.3125
$ws.Range("H132").Value = 36689.895
$ws.Range("I132").Value = 37863.594
$ws.Range("K132").Value = 113590.782
$ws.Range("M132").Value = -111060.782
$ws.Range("H134").Value = 1454.5405
$ws.Range("I134").Value = 1302
$ws.Range("K134").Value = 3906
$ws.Range("M134").Value = -1371

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 305.33334
$ws.Range("I14").Value = 305.33334
$ws.Range("K14").Value = 916.0000200000001
$ws.Range("M14").Value = -743.0000200000001
$ws.Range("H134").Value = 1093
$ws.Range("I134").Value = 1093
$ws.Range("K134").Value = 3279
$ws.Range("M134").Value = 1791
$ws.Range("H136").Value = 2398
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
$ws.Range("H138").Value = 7444.375
$ws.Range("I138").Value = 1753.6666
$ws.Range("J138").Value = 24516.5
$ws.Range("K138").Value = 5260.9998
$ws.Range("L138").Value = 73549.5
$ws.Range("M138").Value = -120.9997999999996
$ws.Range("N138").Value = -83829.5
$ws.Range("H139").Value = 1419.8
$ws.Range("I139").Value = 1299.5
$ws.Range("J139").Value = 1500
$ws.Range("K139").Value = 3898.5
$ws.Range("L139").Value = 4500
$ws.Range("M139").Value = 1241.5
$ws.Range("N139").Value = -14780

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1475.9656
$ws.Range("I97").Value = 922.1579
$ws.Range("J97").Value = 2528.2
$ws.Range("K97").Value = 922.1579
$ws.Range("L97").Value = 2528.2
$ws.Range("M97").Value = -426.1579
$ws.Range("N97").Value = -3520.2
$ws.Range("H102").Value = 26000.04
$ws.Range("I102").Value = 28864
$ws.Range("J102").Value = 4997.6665
$ws.Range("K102").Value = 28864
$ws.Range("L102").Value = 4997.6665
$ws.Range("M102").Value = -27242
$ws.Range("N102").Value = -8241.666499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1735.4722
$ws.Range("I40").Value = 1642.2
$ws.Range("K40").Value = 1642.2
$ws.Range("M40").Value = -1506.2
$ws.Range("H122").Value = 4724.476
$ws.Range("I122").Value = 3773.4375
$ws.Range("J122").Value = 7767.8
$ws.Range("K122").Value = 11320.3125
$ws.Range("L122").Value = 23303.4
$ws.Range("M122").Value = -8870.3125
$ws.Range("N122").Value = -28203.4
$ws.Range("H132").Value = 2692.111
$ws.Range("I132").Value = 1146.6
$ws.Range("K132").Value = 3439.8
$ws.Range("M132").Value = -909.7999999999997
$ws.Range("H136").Value = 2858.4358
$ws.Range("I136").Value = 3068.6667
$ws.Range("K136").Value = 9206.000100000001
$ws.Range("M136").Value = -6656.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2345
$ws.Range("I132").Value = 2631
$ws.Range("J132").Value = 2249.6667
$ws.Range("K132").Value = 7893
$ws.Range("L132").Value = 2631
$ws.Range("M132").Value = -5363
$ws.Range("N132").Value = -11809.0001
